# Auto-generated edit script: updates market-price / profit columns (H-N)
# across several worksheet rows, per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

# ALC row 5 (Leve Item ID 5503)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 146.125
$ws.Range("I5").Value = 133.8
$ws.Range("K5").Value = 133.8
$ws.Range("M5").Value = -18.80000000000001

# ALC row 86 (Leve Item ID 12603)
$ws.Range("H86").Value = 2009.2858
$ws.Range("I86").Value = 1317.7273
$ws.Range("J86").Value = 2770
$ws.Range("K86").Value = 1317.7273
$ws.Range("L86").Value = 2770
$ws.Range("M86").Value = -194.7273
$ws.Range("N86").Value = -5016

# ALC row 89 (Leve Item ID 12603)
$ws.Range("H89").Value = 2009.2858
$ws.Range("I89").Value = 1317.7273
$ws.Range("J89").Value = 2770
$ws.Range("K89").Value = 6588.636500000001
$ws.Range("L89").Value = 13850
$ws.Range("M89").Value = -972.6365000000005
$ws.Range("N89").Value = -25082

# ALC row 127 (Leve Item ID 36114)
$ws.Range("H127").Value = 917.43475
$ws.Range("I127").Value = 456.75
$ws.Range("J127").Value = 1420
$ws.Range("K127").Value = 1370.25
$ws.Range("L127").Value = 4260
$ws.Range("M127").Value = 3589.75
$ws.Range("N127").Value = -14180

# ALC row 131 (Leve Item ID 36108)
$ws.Range("H131").Value = 2504.853
$ws.Range("I131").Value = 1239.4546
$ws.Range("J131").Value = 4824.75
$ws.Range("K131").Value = 3718.3638
$ws.Range("L131").Value = 14474.25
$ws.Range("M131").Value = 1321.6362
$ws.Range("N131").Value = -24554.25

# ALC row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 2099.3594
$ws.Range("I132").Value = 1683.2745
$ws.Range("J132").Value = 3731.6924
$ws.Range("K132").Value = 5049.8235
$ws.Range("L132").Value = 11195.0772
$ws.Range("M132").Value = -2519.8235
$ws.Range("N132").Value = -16255.0772

# ARM row 43 (Leve Item ID 21715)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

# ARM row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 1446.8286
$ws.Range("I61").Value = 1176
$ws.Range("J61").Value = 2123.9
$ws.Range("K61").Value = 1176
$ws.Range("L61").Value = 2123.9
$ws.Range("M61").Value = -964
$ws.Range("N61").Value = -2547.9

# ARM row 104 (Leve Item ID 18672)
$ws.Range("H104").Value = 7333.3335
$ws.Range("J104").Value = 7333.3335
$ws.Range("L104").Value = 7333.3335
$ws.Range("N104").Value = -14321.3335

# ARM row 109 (Leve Item ID 25646)
$ws.Range("H109").Value = 27000
$ws.Range("J109").Value = 27000
$ws.Range("L109").Value = 27000
$ws.Range("N109").Value = -29774

# ARM row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 1446.8286
$ws.Range("I136").Value = 1176
$ws.Range("J136").Value = 2123.9
$ws.Range("K136").Value = 3528
$ws.Range("L136").Value = 6371.700000000001
$ws.Range("M136").Value = -978
$ws.Range("N136").Value = -11471.7

# BSM row 86 (Leve Item ID 12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1963.2
$ws.Range("I86").Value = 1837.9474
$ws.Range("J86").Value = 2179.5454
$ws.Range("K86").Value = 1837.9474
$ws.Range("L86").Value = 2179.5454
$ws.Range("M86").Value = -714.9474
$ws.Range("N86").Value = -4425.5454

# BSM row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 1963.2
$ws.Range("I89").Value = 1837.9474
$ws.Range("J89").Value = 2179.5454
$ws.Range("K89").Value = 9189.737000000001
$ws.Range("L89").Value = 10897.727
$ws.Range("M89").Value = -3573.737000000001
$ws.Range("N89").Value = -22129.727

# BSM row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 1094.9714
$ws.Range("I134").Value = 736.44446
$ws.Range("J134").Value = 2305
$ws.Range("K134").Value = 2209.33338
$ws.Range("L134").Value = 6915
$ws.Range("M134").Value = 325.66662
$ws.Range("N134").Value = -11985

# CRP row 58 (Leve Item ID 44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1057.8667
$ws.Range("I58").Value = 681.7692
$ws.Range("J58").Value = 1345.4706
$ws.Range("K58").Value = 681.7692
$ws.Range("L58").Value = 1345.4706
$ws.Range("M58").Value = -478.7692
$ws.Range("N58").Value = -1751.4706

# CRP row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 2388
$ws.Range("I99").Value = 1650.25
$ws.Range("J99").Value = 3568.4
$ws.Range("K99").Value = 1650.25
$ws.Range("L99").Value = 3568.4
$ws.Range("M99").Value = -152.25
$ws.Range("N99").Value = -6564.4

# CRP row 112 (Leve Item ID 25796)
$ws.Range("H112").Value = 34351
$ws.Range("J112").Value = 34351
$ws.Range("L112").Value = 34351
$ws.Range("N112").Value = -37305

# CRP row 122 (Leve Item ID 36196)
$ws.Range("H122").Value = 1228.5
$ws.Range("I122").Value = 1235.0834
$ws.Range("J122").Value = 1189
$ws.Range("K122").Value = 3705.2502
$ws.Range("L122").Value = 3567
$ws.Range("M122").Value = -1255.2502
$ws.Range("N122").Value = -8467

# CRP row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 2388
$ws.Range("I126").Value = 1650.25
$ws.Range("J126").Value = 3568.4
$ws.Range("K126").Value = 4950.75
$ws.Range("L126").Value = 10705.2
$ws.Range("M126").Value = -2480.75
$ws.Range("N126").Value = -15645.2

# CRP row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 2334.96
$ws.Range("I134").Value = 2165.375
$ws.Range("J134").Value = 2636.4443
$ws.Range("K134").Value = 6496.125
$ws.Range("L134").Value = 7909.3329
$ws.Range("M134").Value = -3961.125
$ws.Range("N134").Value = -12979.3329

# CRP row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 1057.8667
$ws.Range("I136").Value = 681.7692
$ws.Range("J136").Value = 1345.4706
$ws.Range("K136").Value = 2045.3076
$ws.Range("L136").Value = 4036.4118
$ws.Range("M136").Value = 504.6924000000001
$ws.Range("N136").Value = -9136.4118

# CUL row 5 (Leve Item ID 43974)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 872.7273
$ws.Range("J5").Value = 1097.5
$ws.Range("L5").Value = 3292.5
$ws.Range("N5").Value = -3516.5

# CUL row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 856.93335
$ws.Range("I131").Value = 778
$ws.Range("J131").Value = 947.1429000000001
$ws.Range("K131").Value = 2334
$ws.Range("L131").Value = 2841.4287
$ws.Range("M131").Value = 2706
$ws.Range("N131").Value = -12921.4287

# CUL row 135 (Leve Item ID 43974)
$ws.Range("H135").Value = 872.7273
$ws.Range("J135").Value = 1097.5
$ws.Range("L135").Value = 9877.5
$ws.Range("N135").Value = -14947.5

# LTW row 7 (Leve Item ID 36249)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1460.2
$ws.Range("I7").Value = 1189.3
$ws.Range("J7").Value = 2002
$ws.Range("K7").Value = 1189.3
$ws.Range("L7").Value = 2002
$ws.Range("M7").Value = -1077.3
$ws.Range("N7").Value = -2226

# LTW row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 1611.9642
$ws.Range("I40").Value = 1248.4375
$ws.Range("K40").Value = 1248.4375
$ws.Range("M40").Value = -1112.4375

# LTW row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 37990.43
$ws.Range("I122").Value = 45442.26
$ws.Range("J122").Value = 3712
$ws.Range("K122").Value = 136326.78
$ws.Range("L122").Value = 11136
$ws.Range("M122").Value = -133876.78
$ws.Range("N122").Value = -16036

# LTW row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 1460.2
$ws.Range("I126").Value = 1189.3
$ws.Range("J126").Value = 2002
$ws.Range("K126").Value = 3567.9
$ws.Range("L126").Value = 6006
$ws.Range("M126").Value = -1097.9
$ws.Range("N126").Value = -10946

# WVR row 64 (Leve Item ID 11036)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 10000
$ws.Range("J64").Value = 10000
$ws.Range("L64").Value = 10000
$ws.Range("N64").Value = -10496

# WVR row 67 (Leve Item ID 11036)
$ws.Range("H67").Value = 10000
$ws.Range("J67").Value = 10000
$ws.Range("L67").Value = 10000
$ws.Range("N67").Value = -11716

# WVR row 122 (Leve Item ID 36208)
$ws.Range("H122").Value = 9092498
$ws.Range("I122").Value = 11112315
$ws.Range("J122").Value = 3322.5
$ws.Range("K122").Value = 33336945
$ws.Range("L122").Value = 9967.5
$ws.Range("M122").Value = -33334495
$ws.Range("N122").Value = -14867.5

# WVR row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 3158.3333
$ws.Range("I132").Value = 3201.4092
$ws.Range("J132").Value = 2887.5715
$ws.Range("K132").Value = 9604.2276
$ws.Range("L132").Value = 8662.7145
$ws.Range("M132").Value = -7074.2276
$ws.Range("N132").Value = -13722.7145
